# Mise a jour de l'application
# - Adds MD ("M" = Match day) markers to column D for the games that were
#   still missing them (rows 340-351 "N3 J1 VS Montpellier (B)" and
#   398-411 "N3 J2 VS Ales").
# - Appends the 14 player rows (441-454) for the new match
#   "N3 J3 VS OM (B)" played 2025-09-06, each one also tagged MD = "M".
# - Refreshes the sheet selection to match the new bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) New rows' data (row, A..V). String cells are prefixed "S:", numeric
#    cells "N:" so the parsing below can tell them apart without relying on
#    PowerShell's own type inference.
# ---------------------------------------------------------------------------
$newRowsText = @"
441|S:N3 J3 VS OM (B)|N:45906|S:Global|S:M|S:Ilan Ihaddadene|S:center midfield|S:01:18:38|N:9.75|N:2.3|N:7.43|N:1.83|N:0.44|N:0.05|N:0|N:6|N:7.41|N:26.86|N:4.25|N:38|N:2|N:17|N:5
442|S:N3 J3 VS OM (B)|N:45906|S:Global|S:M|S:Amir Etien|S:right forward|S:01:40:20|N:8.52|N:1.57|N:6.94|N:0.78|N:0.45|N:0.35|N:0.01|N:20|N:5.06|N:31.39|N:5.1|N:39|N:13|N:29|N:13
443|S:N3 J3 VS OM (B)|N:45906|S:Global|S:M|S:Naim Ighbane|S:center back|S:01:00:27|N:5.95|N:0.95|N:4.99|N:0.56|N:0.3|N:0.09|N:0|N:6|N:5.88|N:28.65|N:3.74|N:8|N:0|N:15|N:9
444|S:N3 J3 VS OM (B)|N:45906|S:Global|S:M|S:Omar Benyounes|S:center midfield|S:00:23:13|N:2.72|N:0.69|N:2.01|N:0.48|N:0.19|N:0.03|N:0|N:2|N:6.99|N:27.43|N:3.69|N:13|N:0|N:15|N:3
445|S:N3 J3 VS OM (B)|N:45906|S:Global|S:M|S:Ilyes Boughanmi|S:center forward|S:00:23:03|N:2.11|N:0.37|N:1.73|N:0.27|N:0.07|N:0.04|N:0|N:3|N:5.47|N:28.38|N:4.75|N:9|N:2|N:12|N:2
446|S:N3 J3 VS OM (B)|N:45906|S:Global|S:M|S:Mattheo Haon|S:right back|S:01:39:09|N:10.32|N:1.51|N:8.79|N:1|N:0.41|N:0.1|N:0.02|N:10|N:6.23|N:31.69|N:4.78|N:36|N:6|N:37|N:15
447|S:N3 J3 VS OM (B)|N:45906|S:Global|S:M|S:Yoan Zouma|S:center back|S:01:39:39|N:9.2|N:1.62|N:7.56|N:1.07|N:0.46|N:0.12|N:0|N:9|N:5.51|N:29.95|N:4.92|N:31|N:7|N:27|N:9
448|S:N3 J3 VS OM (B)|N:45906|S:Global|S:M|S:Naim Dhib|S:center midfield|S:01:39:30|N:10.31|N:1.77|N:8.51|N:1.39|N:0.38|N:0.02|N:0|N:3|N:6.16|N:27.12|N:4.82|N:52|N:11|N:43|N:14
449|S:N3 J3 VS OM (B)|N:45906|S:Global|S:M|S:Emmanuel Valey|S:left forward|S:01:01:07|N:7.37|N:1.7|N:5.66|N:1.13|N:0.47|N:0.1|N:0.01|N:11|N:7.19|N:30.92|N:4.98|N:26|N:12|N:22|N:13
450|S:N3 J3 VS OM (B)|N:45906|S:Global|S:M|S:Rayane Chayebi|S:center midfield|S:01:40:00|N:9.82|N:1.53|N:8.27|N:1.23|N:0.24|N:0.08|N:0|N:6|N:5.84|N:28.95|N:4.87|N:33|N:11|N:48|N:11
451|S:N3 J3 VS OM (B)|N:45906|S:Global|S:M|S:Levy Ndoutoume|S:left back|S:01:39:19|N:9.72|N:1.76|N:7.93|N:1.11|N:0.55|N:0.13|N:0|N:12|N:5.83|N:29.95|N:4.39|N:42|N:11|N:53|N:14
452|S:N3 J3 VS OM (B)|N:45906|S:Global|S:M|S:Jeremie Laurent|S:left forward|S:00:38:33|N:4.2|N:1.08|N:3.11|N:0.76|N:0.29|N:0.04|N:0|N:4|N:6.55|N:26.78|N:4.47|N:29|N:2|N:24|N:9
453|S:N3 J3 VS OM (B)|N:45906|S:Global|S:M|S:Karahali Souaré|S:right forward|S:00:38:43|N:3.83|N:0.76|N:3.07|N:0.4|N:0.24|N:0.11|N:0.02|N:8|N:5.91|N:31.93|N:5.18|N:20|N:6|N:16|N:8
454|S:N3 J3 VS OM (B)|N:45906|S:Global|S:M|S:Sofiane Belle|S:left forward|S:01:18:47|N:8.13|N:1.56|N:6.55|N:1.09|N:0.36|N:0.12|N:0|N:10|N:6.13|N:30.02|N:4.69|N:24|N:2|N:26|N:9
"@

$newRows = New-Object System.Collections.ArrayList
foreach ($line in ($newRowsText -split "`n")) {
    $line = $line.Trim()
    if ($line.Length -eq 0) { continue }
    $fields = $line -split '\|'
    $rowNum = [int]$fields[0]
    $values = New-Object System.Collections.ArrayList
    for ($i = 1; $i -lt $fields.Count; $i++) {
        $field = $fields[$i]
        $tag = $field.Substring(0, 2)
        $raw = $field.Substring(2)
        if ($tag -eq "S:") {
            $null = $values.Add($raw)
        } else {
            $null = $values.Add([double]$raw)
        }
    }
    $null = $newRows.Add(@{ Row = $rowNum; Values = $values })
}

# ---------------------------------------------------------------------------
# 2) Column G ("Temps joue") for the new rows first - matches the order the
#    workbook's own tool wrote the sheet in (time values land in the shared
#    string table before the "M"/match-name strings).
# ---------------------------------------------------------------------------
foreach ($entry in $newRows) {
    $ws.Cells.Item($entry.Row, 7).Value = $entry.Values[6]
}

# ---------------------------------------------------------------------------
# 3) Column D ("MD") - mark every match-day row, old and new, with "M".
# ---------------------------------------------------------------------------
$matchDayRows = @(340,341,342,343,344,345,346,347,348,349,350,351,398,399,400,401,402,403,404,405,406,407,408,409,410,411,441,442,443,444,445,446,447,448,449,450,451,452,453,454)
foreach ($r in $matchDayRows) {
    $ws.Cells.Item($r, 4).Value = "M"
}

# ---------------------------------------------------------------------------
# 4) Remaining columns (A,B,C,E,F,H..V) for the 14 brand-new rows.
# ---------------------------------------------------------------------------
foreach ($entry in $newRows) {
    $r = $entry.Row
    $vals = $entry.Values
    $ws.Cells.Item($r, 1).Value = $vals[0]    # A - Type (match name)
    $ws.Cells.Item($r, 2).Value = $vals[1]    # B - Date
    $ws.Cells.Item($r, 3).Value = $vals[2]    # C - Periode
    # D already set above
    $ws.Cells.Item($r, 5).Value = $vals[4]    # E - Nom du joueur
    $ws.Cells.Item($r, 6).Value = $vals[5]    # F - Poste
    # G already set above
    for ($col = 8; $col -le 22; $col++) {
        $ws.Cells.Item($r, $col).Value = $vals[$col - 1]
    }
}

# ---------------------------------------------------------------------------
# 5) Copy the date format from an existing date cell onto the new B cells so
#    they keep the workbook's date display (instead of a raw serial number).
# ---------------------------------------------------------------------------
$ws.Range("B2").Copy()
$ws.Range("B441:B454").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 6) Move the active selection to mirror where the sheet was left scrolled
#    to after the update.
# ---------------------------------------------------------------------------
$null = $ws.Range("E458").Select()
